$d = $word.ActiveDocument

# The "{representative_user_code}" placeholder that must be removed lives in
# the last row of the second table (the "NAUDOTOJO KODAS..." row), not in the
# legend table further down the document (which must stay untouched).
$cell = $d.Tables.Item(2).Cell(8, 1)

# Remove the "{representative_user_code}" placeholder entirely, leaving the
# preceding space (which separates it from "{user_code}") intact.
# Wrap = 0 (wdFindStop) and Replace = 1 (wdReplaceOne) keep the operation
# confined to the cell range so the legend text later in the document
# (which also mentions this placeholder name) is left untouched.
$cell.Range.Find.Execute("{representative_user_code}", $true, $false, $false, $false, $false,
                          $true, 0, $false, "", 1)
